# Updated cryptos list on Sat Aug 12 10:55:58 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The whole coin table (B2:E51) is stored as text (prices like "29.413.93"
# or "0.9997" are not real numbers). Temporarily force Text format over the
# range before writing so Excel's COM layer doesn't silently reinterpret
# numeric-looking strings as floating point values; restore the default
# style afterwards so the saved file's formatting is unchanged.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

# --- Rows 2-43: refreshed Price (D) / Volume(1h) (E) figures ---
# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.419.94"
$ws.Range("E2").Value = "  +0.20%  "
# Row 3 - Ethereum
$ws.Range("D3").Value = "1.850.18"
$ws.Range("E3").Value = "  +0.20%  "
# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  +0.09%  "
# Row 5 - BNB
$ws.Range("E5").Value = "  +0.21%  "
# Row 6 - XRP
$ws.Range("D6").Value = "0.6308"
$ws.Range("E6").Value = "  +0.09%  "
# Row 7 - USDC (unchanged)
# Row 8 - Dogecoin
$ws.Range("D8").Value = "0.07717"
$ws.Range("E8").Value = "  +2.31%  "
# Row 9 - Cardano
$ws.Range("D9").Value = "0.2944"
$ws.Range("E9").Value = "  -0.36%  "
# Row 10 - Solana
$ws.Range("D10").Value = "24.52"
$ws.Range("E10").Value = "  +0.44%  "
# Row 11 - TRON
$ws.Range("E11").Value = "  +0.55%  "
# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.852.88"
$ws.Range("E12").Value = "  -0.83%  "
# Row 13 - Polkadot
$ws.Range("D13").Value = "5.031"
$ws.Range("E13").Value = "  +0.74%  "
# Row 14 - ShibaInu
$ws.Range("D14").Value = "0.00001086"
$ws.Range("E14").Value = "  +8.24%  "
# Row 15 - Polygon
$ws.Range("D15").Value = "0.6805"
$ws.Range("E15").Value = "  -0.30%  "
# Row 16 - Litecoin
$ws.Range("D16").Value = "83.75"
$ws.Range("E16").Value = "  +1.07%  "
# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "2.105.32"
$ws.Range("E17").Value = "  -0.71%  "
# Row 18 - Uniswap
$ws.Range("D18").Value = "6.156"
$ws.Range("E18").Value = "  +0.48%  "
# Row 19 - WrappedBTC
$ws.Range("D19").Value = "29.439.23"
$ws.Range("E19").Value = "  +0.16%  "
# Row 20 - BitcoinCash
$ws.Range("D20").Value = "229.60"
$ws.Range("E20").Value = "  +1.04%  "
# Row 21 - Avalanche
$ws.Range("D21").Value = "12.48"
$ws.Range("E21").Value = "  +0.32%  "
# Row 22 - Dai
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.06%  "
# Row 23 - Chainlink
$ws.Range("D23").Value = "7.459"
$ws.Range("E23").Value = "  -1.08%  "
# Row 24 - BinanceUSD
$ws.Range("E24").Value = "  +0.06%  "
# Row 25 - Monero
$ws.Range("D25").Value = "157.57"
$ws.Range("E25").Value = "  +0.32%  "
# Row 26 - Stellar
$ws.Range("E26").Value = "  -0.22%  "
# Row 27 - Cosmos
$ws.Range("D27").Value = "8.377"
$ws.Range("E27").Value = "  +0.30%  "
# Row 28 - EthereumClassic
$ws.Range("D28").Value = "17.69"
$ws.Range("E28").Value = "  +0.26%  "
# Row 29 - PancakeSwap
$ws.Range("D29").Value = "1.469"
$ws.Range("E29").Value = "  +0.25%  "
# Row 30 - Toncoin
$ws.Range("D30").Value = "1.313"
$ws.Range("E30").Value = "  +4.89%  "
# Row 31 - Hedera
$ws.Range("D31").Value = "0.05745"
$ws.Range("E31").Value = "  +1.19%  "
# Row 32 - Filecoin
$ws.Range("D32").Value = "4.114"
$ws.Range("E32").Value = "  -0.12%  "
# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "4.055"
$ws.Range("E33").Value = "  +0.96%  "
# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +0.49%  "
# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  +0.58%  "
# Row 36 - ImmutableX
$ws.Range("E36").Value = "  -0.43%  "
# Row 37 - HuobiToken
$ws.Range("E37").Value = "  -0.34%  "
# Row 38 - MXToken
$ws.Range("D38").Value = "2.778"
$ws.Range("E38").Value = "  -0.03%  "
# Row 39 - Maker
$ws.Range("D39").Value = "1.229.31"
$ws.Range("E39").Value = "  -2.33%  "
# Row 40 - VeChain
$ws.Range("E40").Value = "  -0.90%  "
# Row 41 - FraxShare
$ws.Range("D41").Value = "6.474"
$ws.Range("E41").Value = "  +4.15%  "
# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "0.9122"
$ws.Range("E42").Value = "  +0.07%  "
# Row 43 - PaxDollar
$ws.Range("E43").Value = "  +0.06%  "

# --- Rows 44-51: RocketPoolETH enters the ranking, pushing the rest down
#     by one slot; Algorand (previously row 51) drops off the list. ---
# Row 44 - RocketPoolETH (new entry)
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.014.07"
$ws.Range("E44").Value = "  -0.73%  "
# Row 45 - Quant
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "101.76"
$ws.Range("E45").Value = "  +0.56%  "
# Row 46 - Aave
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "66.36"
$ws.Range("E46").Value = "  +0.29%  "
# Row 47 - BabyDogeCoin
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.00000000122"
$ws.Range("E47").Value = "  +3.24%  "
# Row 48 - Aptos
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.154"
$ws.Range("E48").Value = "  +1.30%  "
# Row 49 - TheSandbox
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "0.4019"
$ws.Range("E49").Value = "  -0.45%  "
# Row 50 - EnergySwap
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.020"
$ws.Range("E50").Value = "  -1.00%  "
# Row 51 - RenderToken (Algorand, formerly here, is dropped)
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "1.689"
$ws.Range("E51").Value = "  +0.56%  "

# Restore default styling over the data range (drops the temporary text
# NumberFormat so the saved cells carry the same style they started with).
$dataRange.Style = "Normal"
